$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 229, pushing the existing rows 229..259 down to 230..260
$ws.Rows(229).Insert()

# Populate the newly inserted row 229 with the new Jengibre price record
$ws.Cells.Item(229,1).Value = 10
$ws.Cells.Item(229,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(229,3).Value = "La Araucanía"
$ws.Cells.Item(229,4).Value = 44984
$ws.Cells.Item(229,5).Value = 9
$ws.Cells.Item(229,6).Value = 100114007
$ws.Cells.Item(229,7).Value = "Jengibre"
$ws.Cells.Item(229,8).Value = "Sin especificar"
$ws.Cells.Item(229,9).Value = "Primera"
$ws.Cells.Item(229,10).Value = 100
$ws.Cells.Item(229,11).Value = 28000
$ws.Cells.Item(229,12).Value = 28000
$ws.Cells.Item(229,13).Value = 28000
$ws.Cells.Item(229,14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(229,15).Value = "Perú"
$ws.Cells.Item(229,16).Value = 2154
$ws.Cells.Item(229,17).Value = 13
$ws.Cells.Item(229,18).Value = "Hortaliza"
